$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.359.97'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '3.503.82'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'586.27"
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = "'135.67"
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('D7').Value = '3.504.21'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  -3.45%  '
$ws.Range('D13').Value = '4.099.53'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '3.502.78'
$ws.Range('D17').Value = '64.350.42'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = "'25.43"
$ws.Range('E18').Value = '  -8.78%  '
$ws.Range('D19').Value = "'9.79"
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').Value = "'5.59"
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('D22').Value = "'383.93"
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('D23').Value = "'0.570"
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('D24').Value = '3.641.81'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = "'73.98"
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('E28').Value = '  +3.64%  '
$ws.Range('D29').Value = "'1.56"
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').Value = "'7.53"
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').Value = "'0.999"
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = "'8.29"
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').Value = '3.524.48'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = "'23.56"
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').Value = "'5.32"
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').Value = "'163.78"
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').Value = "'0.0787"
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('D43').Value = "'0.807"
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = "'25.85"
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = "'41.92"
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').Value = "'4.41"
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = '2.473.28'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('E51').Value = '  -1.91%  '
